$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in B1 from "Azoto (g)" to "Azoto"
$ws.Range("B1").Value = "Azoto"

# Update the active cell selection to B1
$ws.Range("B1").Select()
